$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dMap = @{
    2 = "28.674.06"
    3 = "1.803.93"
    5 = "231.49"
    6 = "0.5964"
    8 = "0.2782"
    9 = "0.06849"
    10 = "23.39"
    11 = "0.07523"
    12 = "1.798.11"
    13 = "4.707"
    14 = "0.6268"
    15 = "2.048.09"
    16 = "0.000009193"
    17 = "75.35"
    18 = "28.648.77"
    19 = "5.462"
    21 = "210.48"
    22 = "11.43"
    23 = "6.849"
    25 = "154.46"
    26 = "7.843"
    27 = "0.1277"
    29 = "1.442"
    30 = "0.06285"
    32 = "3.766"
    33 = "3.734"
    34 = "1.712"
    35 = "1.053"
    36 = "0.6364"
    37 = "2.493"
    38 = "2.717"
    39 = "0.01707"
    40 = "6.390"
    41 = "1.136.43"
    42 = "0.8666"
    44 = "100.72"
    45 = "1.962.74"
    46 = "60.55"
    47 = "0.00000000111"
    48 = "1.581"
    49 = "8.322"
    51 = "0.05440"
}

$eMap = @{
    2 = "  -1.90%  "
    3 = "  -1.45%  "
    4 = "  +0.21%  "
    5 = "  -2.17%  "
    6 = "  -1.72%  "
    7 = "  +0.15%  "
    8 = "  -1.31%  "
    9 = "  -3.72%  "
    10 = "  -2.36%  "
    11 = "  -1.93%  "
    12 = "  -1.75%  "
    13 = "  -2.61%  "
    14 = "  -1.69%  "
    15 = "  -1.49%  "
    16 = "  -9.09%  "
    17 = "  -5.24%  "
    18 = "  -1.95%  "
    19 = "  -7.58%  "
    20 = "  +0.18%  "
    21 = "  -8.07%  "
    22 = "  -3.43%  "
    23 = "  -2.61%  "
    24 = "  +0.26%  "
    25 = "  +0.03%  "
    26 = "  -3.06%  "
    27 = "  -1.32%  "
    28 = "  -1.65%  "
    29 = "  -3.15%  "
    30 = "  -2.99%  "
    31 = "  -2.70%  "
    32 = "  -1.70%  "
    33 = "  -2.53%  "
    34 = "  -1.91%  "
    35 = "  -6.74%  "
    36 = "  -2.66%  "
    37 = "  -2.56%  "
    38 = "  -1.50%  "
    39 = "  -2.50%  "
    40 = "  -1.95%  "
    41 = "  -7.00%  "
    42 = "  -7.03%  "
    43 = "  +0.22%  "
    44 = "  -0.41%  "
    45 = "  -0.83%  "
    46 = "  -4.57%  "
    47 = "  -5.26%  "
    48 = "  -1.91%  "
    49 = "  -2.62%  "
    50 = "  -1.58%  "
    51 = "  -1.86%  "
}

foreach ($row in $dMap.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dMap[$row]
}

foreach ($row in $eMap.Keys) {
    $ws.Range("E$row").Value = $eMap[$row]
}
